$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Case ID", "Study Code", "Study Type", "Breed", "Diagnosis", "Stage of Disease", "Age", "Sex", "Neutered Status")
$values  = @("NCATS-COP01CCB010072", "NCATS-COP01", "Transcriptomics", "Akita", "Bone sarcomas :: Osteosarcoma (appendicular)", "", "10", "Female", "Yes")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(2, $col)
    $text = $values[$i]
    if ($text -eq "") {
        # Force an empty-string text cell instead of a blank/deleted cell.
        $cell.Value = "'"
    } elseif ($text -match '^[0-9]+$') {
        # Force text storage so numeric-looking strings (e.g. "10") stay text
        # instead of being auto-converted to a number.
        $cell.NumberFormat = "@"
        $cell.Value = $text
    } else {
        $cell.Value = $text
    }
}
